$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain-text values such as "34.561.48" or
# "1.00" (thousand-separator/trailing-zero formatted numbers kept as text).
# Assigning a numeric-looking string straight to .Value lets Excel coerce it
# to a real number, which would lose the trailing zero / text formatting.
# Temporarily force the whole Price column to Text format so every new value
# we write lands as a string, then restore the default ("Normal") style so
# we do not leave a stray number-format override behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "34.561.48"
$ws.Range("E2").Value = "  +13.79%  "

# Row 3
$ws.Range("D3").Value = "1.821.66"
$ws.Range("E3").Value = "  +8.11%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.25%  "

# Row 5
$ws.Range("D5").Value = "233.19"
$ws.Range("E5").Value = "  +5.41%  "

# Row 6
$ws.Range("D6").Value = "0.549"
$ws.Range("E6").Value = "  +5.11%  "

# Row 7
$ws.Range("E7").Value = "  +0.25%  "

# Row 8
$ws.Range("D8").Value = "31.69"
$ws.Range("E8").Value = "  +4.47%  "

# Row 9
$ws.Range("D9").Value = "46.09"
$ws.Range("E9").Value = "  +4.42%  "

# Row 10
$ws.Range("D10").Value = "0.285"
$ws.Range("E10").Value = "  +7.60%  "

# Row 11
$ws.Range("E11").Value = "  +9.64%  "

# Row 12
$ws.Range("D12").Value = "0.0932"
$ws.Range("E12").Value = "  +3.41%  "

# Row 13
$ws.Range("D13").Value = "2.086.26"
$ws.Range("E13").Value = "  +8.30%  "

# Row 14
$ws.Range("D14").Value = "1.850.08"
$ws.Range("E14").Value = "  +9.59%  "

# Row 15
$ws.Range("D15").Value = "0.647"
$ws.Range("E15").Value = "  +4.56%  "

# Row 16
$ws.Range("D16").Value = "34.511.54"
$ws.Range("E16").Value = "  +13.68%  "

# Row 17
$ws.Range("D17").Value = "10.29"
$ws.Range("E17").Value = "  -4.15%  "

# Row 18
$ws.Range("D18").Value = "4.34"
$ws.Range("E18").Value = "  +8.18%  "

# Row 19
$ws.Range("D19").Value = "71.16"
$ws.Range("E19").Value = "  +8.00%  "

# Row 20
$ws.Range("D20").Value = "261.95"
$ws.Range("E20").Value = "  +6.18%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0759"
$ws.Range("E21").Value = "  +5.25%  "

# Row 22
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.20%  "

# Row 23
$ws.Range("D23").Value = "10.54"
$ws.Range("E23").Value = "  +3.40%  "

# Row 24
$ws.Range("E24").Value = "  +2.69%  "

# Row 25
$ws.Range("E25").Value = "  -0.43%  "

# Row 26
$ws.Range("D26").Value = "162.25"
$ws.Range("E26").Value = "  +2.20%  "

# Row 27
$ws.Range("D27").Value = "16.95"
$ws.Range("E27").Value = "  +6.77%  "

# Row 28
$ws.Range("E28").Value = "  +5.13%  "

# Row 29
$ws.Range("E29").Value = "  +6.00%  "

# Row 30
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.33%  "

# Row 31
$ws.Range("D31").Value = "3.84"
$ws.Range("E31").Value = "  +9.39%  "

# Row 32
$ws.Range("E32").Value = "  +3.00%  "

# Row 33
$ws.Range("E33").Value = "  +6.93%  "

# Row 34
$ws.Range("D34").Value = "3.58"
$ws.Range("E34").Value = "  +8.09%  "

# Row 35
$ws.Range("D35").Value = "1.590.94"
$ws.Range("E35").Value = "  +5.34%  "

# Row 36
$ws.Range("E36").Value = "  +6.10%  "

# Row 37
$ws.Range("E37").Value = "  +3.51%  "

# Row 38
$ws.Range("D38").Value = "86.06"
$ws.Range("E38").Value = "  +8.85%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0190"
$ws.Range("E39").Value = "  +5.37%  "

# Row 40
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "0.632"
$ws.Range("E40").Value = "  +7.93%  "

# Row 41
$ws.Range("D41").Value = "2.81"
$ws.Range("E41").Value = "  +2.89%  "

# Row 42
$ws.Range("D42").Value = "2.37"
$ws.Range("E42").Value = "  +2.19%  "

# Row 43
$ws.Range("E43").Value = "  +8.16%  "

# Row 44
$ws.Range("E44").Value = "  +6.78%  "

# Row 45
$ws.Range("D45").Value = "0.0525"
$ws.Range("E45").Value = "  +3.93%  "

# Row 46
$ws.Range("E46").Value = "  +6.42%  "

# Row 47
$ws.Range("D47").Value = "1.976.31"
$ws.Range("E47").Value = "  +8.45%  "

# Row 48
$ws.Range("D48").Value = "53.81"
$ws.Range("E48").Value = "  +3.82%  "

# Row 49
$ws.Range("E49").Value = "  +5.47%  "

# Row 50
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.15%  "

# Row 51
$ws.Range("D51").Value = "0.0₆0124"
$ws.Range("E51").Value = "  +8.92%  "

# Restore the default style on the Price column now that the text values
# are safely written (keeps styles.xml identical to the original aside from
# the cell content changes above).
$priceRange.Style = "Normal"
